{"js": "const replacements = [\n  [\"914\u00f74=228, 2\", \"688\u00f74=172, 0\"],\n  [\"785\u00f78=98, 1\", \"706\u00f75=141, 1\"],\n  [\"451\u00f78=56, 3\", \"392\u00f72=196, 0\"],\n  [\"704\u00f74=176, 0\", \"501\u00f75=100, 1\"],\n  [\"180\u00f74=45, 0\", \"402\u00f74=100, 2\"],\n  [\"125\u00f72=62, 1\", \"638\u00f76=106, 2\"],\n  [\"748\u00f78=93, 4\", \"679\u00f72=339, 1\"],\n  [\"685\u00f78=85, 5\", \"922\u00f79=102, 4\"],\n  [\"762\u00f74=190, 2\", \"630\u00f75=126, 0\"],\n  [\"703\u00f78=87, 7\", \"395\u00f78=49, 3\"],\n  [\"398\u00f76=66, 2\", \"123\u00f74=30, 3\"],\n  [\"714\u00f74=178, 2\", \"590\u00f78=73, 6\"],\n  [\"312\u00f78=39, 0\", \"890\u00f73=296, 2\"],\n  [\"686\u00f76=114, 2\", \"525\u00f73=175, 0\"],\n  [\"452\u00f79=50, 2\", \"832\u00f78=104, 0\"],\n  [\"199\u00f75=39, 4\", \"233\u00f75=46, 3\"],\n  [\"701\u00f74=175, 1\", \"405\u00f75=81, 0\"],\n  [\"689\u00f72=344, 1\", \"404\u00f76=67, 2\"],\n  [\"920\u00f78=115, 0\", \"700\u00f76=116, 4\"],\n  [\"337\u00f76=56, 1\", \"169\u00f73=56, 1\"],\n  [\"973\u00f77=139, 0\", \"104\u00f75=20, 4\"],\n  [\"969\u00f74=242, 1\", \"317\u00f75=63, 2\"],\n  [\"292\u00f72=146, 0\", \"539\u00f74=134, 3\"],\n  [\"948\u00f72=474, 0\", \"797\u00f75=159, 2\"],\n  [\"198\u00f79=22, 0\", \"887\u00f79=98, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"914\u00f74=228, 2\", \"688\u00f74=172, 0\")\n    ,@(\"785\u00f78=98, 1\", \"706\u00f75=141, 1\")\n    ,@(\"451\u00f78=56, 3\", \"392\u00f72=196, 0\")\n    ,@(\"704\u00f74=176, 0\", \"501\u00f75=100, 1\")\n    ,@(\"180\u00f74=45, 0\", \"402\u00f74=100, 2\")\n    ,@(\"125\u00f72=62, 1\", \"638\u00f76=106, 2\")\n    ,@(\"748\u00f78=93, 4\", \"679\u00f72=339, 1\")\n    ,@(\"685\u00f78=85, 5\", \"922\u00f79=102, 4\")\n    ,@(\"762\u00f74=190, 2\", \"630\u00f75=126, 0\")\n    ,@(\"703\u00f78=87, 7\", \"395\u00f78=49, 3\")\n    ,@(\"398\u00f76=66, 2\", \"123\u00f74=30, 3\")\n    ,@(\"714\u00f74=178, 2\", \"590\u00f78=73, 6\")\n    ,@(\"312\u00f78=39, 0\", \"890\u00f73=296, 2\")\n    ,@(\"686\u00f76=114, 2\", \"525\u00f73=175, 0\")\n    ,@(\"452\u00f79=50, 2\", \"832\u00f78=104, 0\")\n    ,@(\"199\u00f75=39, 4\", \"233\u00f75=46, 3\")\n    ,@(\"701\u00f74=175, 1\", \"405\u00f75=81, 0\")\n    ,@(\"689\u00f72=344, 1\", \"404\u00f76=67, 2\")\n    ,@(\"920\u00f78=115, 0\", \"700\u00f76=116, 4\")\n    ,@(\"337\u00f76=56, 1\", \"169\u00f73=56, 1\")\n    ,@(\"973\u00f77=139, 0\", \"104\u00f75=20, 4\")\n    ,@(\"969\u00f74=242, 1\", \"317\u00f75=63, 2\")\n    ,@(\"292\u00f72=146, 0\", \"539\u00f74=134, 3\")\n    ,@(\"948\u00f72=474, 0\", \"797\u00f75=159, 2\")\n    ,@(\"198\u00f79=22, 0\", \"887\u00f79=98, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
